$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: existing "Desk" / "Main Desktop" check-out labels replace the old placeholder text ---
$ws.Range("I15").Value = " Desk"
$ws.Range("J15").Value = " Main Desktop"

# Helper pattern: the A column holds numeric-looking ids ("119","120","121") that must stay
# text (matching the existing "100".."118" ids), so each is written while the cell is
# temporarily formatted as Text, then the format is cleared again so no stray style sticks
# around (keeps the cell on the default "General" style like every other cell in the sheet).

# --- Row 19: new check in/out sample row (beans1 / beans2) ---
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "119"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = " beans1"
$ws.Range("C19").Value = " beans1"
$ws.Range("D19").Value = " beans1"
$ws.Range("E19").Value = " beans1"
$ws.Range("F19").Value = " beans1"
$ws.Range("G19").Value = " beans1"
$ws.Range("H19").Value = " beans1"
$ws.Range("I19").Value = " beans2"
$ws.Range("J19").Value = " beans2"

# --- Row 20: new check in/out sample row (beans3_1 / Storage_1 / N/A_1) ---
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "120"
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = " beans3_1"
$ws.Range("C20").Value = " beans3_1"
$ws.Range("D20").Value = " beans3_1"
$ws.Range("E20").Value = " beans3_1"
$ws.Range("F20").Value = " beans3_1"
$ws.Range("G20").Value = " beans3_1"
$ws.Range("H20").Value = " beans3_1"
$ws.Range("I20").Value = " Storage_1"
$ws.Range("J20").Value = " N/A_1"

# --- Row 21: new check in/out sample row (beans3_1 / Storage / N/A) ---
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "121"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = " beans3_1"
$ws.Range("C21").Value = " beans3_1"
$ws.Range("D21").Value = " beans3_1"
$ws.Range("E21").Value = " beans3_1"
$ws.Range("F21").Value = " beans3_1"
$ws.Range("G21").Value = " beans3_1"
$ws.Range("H21").Value = " beans3_1"
$ws.Range("I21").Value = " Storage"
$ws.Range("J21").Value = " N/A"
